$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Adam9"
$row2[0,2] = "Itga3"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 10.36733066666667
$row2[0,7] = 31.101992
$row2[0,8] = 0.1169328841728879
$row2[0,9] = 0.1169328841728879
$row2[0,10] = 2
$row2[0,11] = 0.6666666666666666
$row2[0,12] = 5.970993
$row2[0,13] = 17.912979
$row2[0,14] = 0.6157237531330177
$row2[0,15] = 0.6157237531330177
$row2[0,16] = 61.90325883935199
$row2[0,17] = 557.129329554168
$row2[0,18] = 0.071998354307599
$row2[0,19] = 0.071998354307599
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Adam9"
$row3[0,2] = "Itga3"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 10.36733066666667
$row3[0,7] = 31.101992
$row3[0,8] = 0.1169328841728879
$row3[0,9] = 0.1169328841728879
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 1.061748
$row3[0,13] = 3.185244
$row3[0,14] = 0.10948655666511
$row3[0,15] = 0.10948655666511
$row3[0,16] = 11.007492600672
$row3[0,17] = 99.06743340604798
$row3[0,18] = 0.01280257884900964
$row3[0,19] = 0.01280257884900964
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Adam9"
$row4[0,2] = "Itga3"
$row4[0,3] = "M2"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 10.36733066666667
$row4[0,7] = 31.101992
$row4[0,8] = 0.1169328841728879
$row4[0,9] = 0.1169328841728879
$row4[0,10] = 1
$row4[0,11] = 0.3333333333333333
$row4[0,12] = 0.02952066666666667
$row4[0,13] = 0.088562
$row4[0,14] = 0.003044146203987976
$row4[0,15] = 0.003044146203987975
$row4[0,16] = 0.3060505128337778
$row4[0,17] = 2.754454615504
$row4[0,18] = 0.0003559607954762624
$row4[0,19] = 0.0003559607954762623
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Adam9"
$row5[0,2] = "Itga3"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 10.36733066666667
$row5[0,7] = 31.101992
$row5[0,8] = 0.1169328841728879
$row5[0,9] = 0.1169328841728879
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 2.635257666666666
$row5[0,13] = 7.905773
$row5[0,14] = 0.2717455439978843
$row5[0,15] = 0.2717455439978843
$row5[0,16] = 27.32058762220177
$row5[0,17] = 245.885288599816
$row5[0,18] = 0.03177599022080302
$row5[0,19] = 0.03177599022080302
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Adam9"
$row6[0,2] = "Itga3"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 37.91490933333333
$row6[0,7] = 113.744728
$row6[0,8] = 0.4276413904453658
$row6[0,9] = 0.4276413904453659
$row6[0,10] = 2
$row6[0,11] = 0.6666666666666666
$row6[0,12] = 5.970993
$row6[0,13] = 17.912979
$row6[0,14] = 0.6157237531330177
$row6[0,15] = 0.6157237531330177
$row6[0,16] = 226.389658224968
$row6[0,17] = 2037.506924024712
$row6[0,18] = 0.2633089619200428
$row6[0,19] = 0.2633089619200429
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Adam9"
$row7[0,2] = "Itga3"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 37.91490933333333
$row7[0,7] = 113.744728
$row7[0,8] = 0.4276413904453658
$row7[0,9] = 0.4276413904453659
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 1.061748
$row7[0,13] = 3.185244
$row7[0,14] = 0.10948655666511
$row7[0,15] = 0.10948655666511
$row7[0,16] = 40.256079154848
$row7[0,17] = 362.304712393632
$row7[0,18] = 0.04682098332734298
$row7[0,19] = 0.04682098332734299
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "FAPs"
$row8[0,1] = "Adam9"
$row8[0,2] = "Itga3"
$row8[0,3] = "M2"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 37.91490933333333
$row8[0,7] = 113.744728
$row8[0,8] = 0.4276413904453658
$row8[0,9] = 0.4276413904453659
$row8[0,10] = 1
$row8[0,11] = 0.3333333333333333
$row8[0,12] = 0.02952066666666667
$row8[0,13] = 0.088562
$row8[0,14] = 0.003044146203987976
$row8[0,15] = 0.003044146203987975
$row8[0,16] = 1.119273400126222
$row8[0,17] = 10.073460601136
$row8[0,18] = 0.0013018029153924
$row8[0,19] = 0.0013018029153924
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "FAPs"
$row9[0,1] = "Adam9"
$row9[0,2] = "Itga3"
$row9[0,3] = "sCs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 37.91490933333333
$row9[0,7] = 113.744728
$row9[0,8] = 0.4276413904453658
$row9[0,9] = 0.4276413904453659
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 2.635257666666666
$row9[0,13] = 7.905773
$row9[0,14] = 0.2717455439978843
$row9[0,15] = 0.2717455439978843
$row9[0,16] = 99.91555550163822
$row9[0,17] = 899.2399995147441
$row9[0,18] = 0.1162096422825876
$row9[0,19] = 0.1162096422825876
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "M2"
$row10[0,1] = "Adam9"
$row10[0,2] = "Itga3"
$row10[0,3] = "ECs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 26.72147866666667
$row10[0,7] = 80.164436
$row10[0,8] = 0.3013909433702152
$row10[0,9] = 0.3013909433702153
$row10[0,10] = 2
$row10[0,11] = 0.6666666666666666
$row10[0,12] = 5.970993
$row10[0,13] = 17.912979
$row10[0,14] = 0.6157237531330177
$row10[0,15] = 0.6157237531330177
$row10[0,16] = 159.553762068316
$row10[0,17] = 1435.983858614844
$row10[0,18] = 0.1855735628122097
$row10[0,19] = 0.1855735628122098
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = "M2"
$row11[0,1] = "Adam9"
$row11[0,2] = "Itga3"
$row11[0,3] = "FAPs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 26.72147866666667
$row11[0,7] = 80.164436
$row11[0,8] = 0.3013909433702152
$row11[0,9] = 0.3013909433702153
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 1.061748
$row11[0,13] = 3.185244
$row11[0,14] = 0.10948655666511
$row11[0,15] = 0.10948655666511
$row11[0,16] = 28.371476531376
$row11[0,17] = 255.343288782384
$row11[0,18] = 0.03299825659965404
$row11[0,19] = 0.03299825659965404
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = "M2"
$row12[0,1] = "Adam9"
$row12[0,2] = "Itga3"
$row12[0,3] = "M2"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 26.72147866666667
$row12[0,7] = 80.164436
$row12[0,8] = 0.3013909433702152
$row12[0,9] = 0.3013909433702153
$row12[0,10] = 1
$row12[0,11] = 0.3333333333333333
$row12[0,12] = 0.02952066666666667
$row12[0,13] = 0.088562
$row12[0,14] = 0.003044146203987976
$row12[0,15] = 0.003044146203987975
$row12[0,16] = 0.7888358645591111
$row12[0,17] = 7.099522781031999
$row12[0,18] = 0.0009174780961767956
$row12[0,19] = 0.0009174780961767956
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = "M2"
$row13[0,1] = "Adam9"
$row13[0,2] = "Itga3"
$row13[0,3] = "sCs"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 26.72147866666667
$row13[0,7] = 80.164436
$row13[0,8] = 0.3013909433702152
$row13[0,9] = 0.3013909433702153
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 2.635257666666666
$row13[0,13] = 7.905773
$row13[0,14] = 0.2717455439978843
$row13[0,15] = 0.2717455439978843
$row13[0,16] = 70.4179815210031
$row13[0,17] = 633.761833689028
$row13[0,18] = 0.08190164586217467
$row13[0,19] = 0.08190164586217467
$ws.Range("A13:T13").Value = $row13

$row14 = New-Object 'object[,]' 1,20
$row14[0,0] = "sCs"
$row14[0,1] = "Adam9"
$row14[0,2] = "Itga3"
$row14[0,3] = "ECs"
$row14[0,4] = 3
$row14[0,5] = 1
$row14[0,6] = 13.65680433333333
$row14[0,7] = 40.970413
$row14[0,8] = 0.154034782011531
$row14[0,9] = 0.154034782011531
$row14[0,10] = 2
$row14[0,11] = 0.6666666666666666
$row14[0,12] = 5.970993
$row14[0,13] = 17.912979
$row14[0,14] = 0.6157237531330177
$row14[0,15] = 0.6157237531330177
$row14[0,16] = 81.544683076703
$row14[0,17] = 733.902147690327
$row14[0,18] = 0.09484287409316612
$row14[0,19] = 0.09484287409316614
$ws.Range("A14:T14").Value = $row14

$row15 = New-Object 'object[,]' 1,20
$row15[0,0] = "sCs"
$row15[0,1] = "Adam9"
$row15[0,2] = "Itga3"
$row15[0,3] = "FAPs"
$row15[0,4] = 3
$row15[0,5] = 1
$row15[0,6] = 13.65680433333333
$row15[0,7] = 40.970413
$row15[0,8] = 0.154034782011531
$row15[0,9] = 0.154034782011531
$row15[0,10] = 3
$row15[0,11] = 1
$row15[0,12] = 1.061748
$row15[0,13] = 3.185244
$row15[0,14] = 0.10948655666511
$row15[0,15] = 0.10948655666511
$row15[0,16] = 14.500084687308
$row15[0,17] = 130.500762185772
$row15[0,18] = 0.01686473788910336
$row15[0,19] = 0.01686473788910336
$ws.Range("A15:T15").Value = $row15

$row16 = New-Object 'object[,]' 1,20
$row16[0,0] = "sCs"
$row16[0,1] = "Adam9"
$row16[0,2] = "Itga3"
$row16[0,3] = "M2"
$row16[0,4] = 3
$row16[0,5] = 1
$row16[0,6] = 13.65680433333333
$row16[0,7] = 40.970413
$row16[0,8] = 0.154034782011531
$row16[0,9] = 0.154034782011531
$row16[0,10] = 1
$row16[0,11] = 0.3333333333333333
$row16[0,12] = 0.02952066666666667
$row16[0,13] = 0.088562
$row16[0,14] = 0.003044146203987976
$row16[0,15] = 0.003044146203987975
$row16[0,16] = 0.4031579684562223
$row16[0,17] = 3.628421716106
$row16[0,18] = 0.0004689043969425174
$row16[0,19] = 0.0004689043969425175
$ws.Range("A16:T16").Value = $row16

$row17 = New-Object 'object[,]' 1,20
$row17[0,0] = "sCs"
$row17[0,1] = "Adam9"
$row17[0,2] = "Itga3"
$row17[0,3] = "sCs"
$row17[0,4] = 3
$row17[0,5] = 1
$row17[0,6] = 13.65680433333333
$row17[0,7] = 40.970413
$row17[0,8] = 0.154034782011531
$row17[0,9] = 0.154034782011531
$row17[0,10] = 3
$row17[0,11] = 1
$row17[0,12] = 2.635257666666666
$row17[0,13] = 7.905773
$row17[0,14] = 0.2717455439978843
$row17[0,15] = 0.2717455439978843
$row17[0,16] = 35.98919832158322
$row17[0,17] = 323.902784894249
$row17[0,18] = 0.04185826563231902
$row17[0,19] = 0.04185826563231902
$ws.Range("A17:T17").Value = $row17
